# Update the "取得日時" (retrieved datetime) column on the "ランサーズ" sheet.
# All data rows (2-13) are re-stamped with the latest scrape time,
# 2025-10-25 12:42:53, replacing the previous 2025-10-25 12:32:59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-25 12:42:53"

$ws.Range("A2:A13").Value = $newTimestamp
